# Applies the "insertion sort" benchmark re-run values, fixes the
# cell D3 style (it was re-entered without the previous numeric style),
# and updates sheet selections / window state to match the new save.

$wb = $excel.ActiveWorkbook

# --- Update "insertion sort" worksheet data (new benchmark timings) ---
$wsInsertion = $wb.Worksheets.Item("insertion sort")

$newValues = @(
    0.028821706771850499,
    0.16448330879211401,
    0.36543011665344199,
    0.69101428985595703,
    1.01195263862609,
    1.45066595077514,
    2.0582830905914302,
    2.85009288787841,
    3.4401211738586399,
    4.1274387836456299
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 3 + $i
    $wsInsertion.Cells.Item($row, 4).Value = $newValues[$i]
}

# D3 lost its numeric cell style when the value was retyped, reset it to
# the default "Normal" style while leaving D4:D12 formatting untouched.
$wsInsertion.Cells.Item(3, 4).Style = "Normal"

# Selection/active cell on this sheet moved to F11 and the tab became the
# active one.
$wsInsertion.Activate()
$wsInsertion.Range("F11").Select()

# --- "python sort" worksheet is no longer the active tab ---
$wsPython = $wb.Worksheets.Item("python sort")
$wsPython.Range("E19").Select()

# --- Restore the window focus on "insertion sort" (activeTab=2) and
#     update the saved window geometry. ---
$wsInsertion.Activate()
$wsInsertion.Range("F11").Select()

$excel.ActiveWindow.WindowState = -4143  # xlNormal
$excel.Left = 2304
$excel.Top = 2304
$excel.Width = 17280
$excel.Height = 8964
